# Validacion de nombres y colores al grafico
# Adds a final "FIN DEL PROGRAMA" marker row to the data range so that
# downstream validation/chart logic can detect the end of the program list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Tabla1" table (APELLIDO/NOMBRE/NOTA1/NOTA2, A1:D9) grows by one row.
# Use ListRows.Add() so the table definition (ref + autoFilter) expands
# from A1:D9 to A1:D10, exactly like using Excel's UI to extend the table.
$tabla1 = $ws.ListObjects.Item("Tabla1")
$tabla1.ListRows.Add() | Out-Null

# New row 10 data: A10 = 0 (sentinel), G10 = "FIN DEL PROGRAMA" (new shared string)
$ws.Range("A10").Value = 0
$ws.Range("G10").Value = "FIN DEL PROGRAMA"

# Update the active selection to reflect where the user ended up afterwards.
$ws.Range("C14").Select()
